$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 16,16
$arr[0,0] = [double]"-1.036232854677933e-22"
$arr[0,1] = [double]"-1.681332379119761e-32"
$arr[0,2] = [double]"-8.166549349893839e-23"
$arr[0,3] = [double]"8.166690288005269e-23"
$arr[0,4] = [double]"3.111256778370345e-33"
$arr[0,5] = [double]"-4.617119692812142e-33"
$arr[0,6] = [double]"4.6171028718429e-33"
$arr[0,7] = [double]"2.71564472260426e-22"
$arr[0,8] = [double]"-4.658119932132633e-22"
$arr[0,9] = [double]"8.042144018329229e-23"
$arr[0,10] = [double]"3.996898153757778e-22"
$arr[0,11] = [double]"9.657171269114655e-12"
$arr[0,12] = [double]"1.016544344117341e-11"
$arr[0,13] = [double]"1.073019029901643e-11"
$arr[0,14] = [double]"1.136137796366431e-11"
$arr[0,15] = [double]"1"
$arr[1,0] = [double]"-1.07301896317416e-11"
$arr[1,1] = [double]"-1.93495569359993e-21"
$arr[1,2] = [double]"7.066125231047504e-18"
$arr[1,3] = [double]"-7.067056197319154e-18"
$arr[1,4] = [double]"4.002330981565128e-22"
$arr[1,5] = [double]"-1.777411577244597e-22"
$arr[1,6] = [double]"-1.596100714375952e-22"
$arr[1,7] = [double]"2.146041825515863e-11"
$arr[1,8] = [double]"-2.154232952819525e-11"
$arr[1,9] = [double]"2.234731249473008e-11"
$arr[1,10] = [double]"1.631888509184236e-17"
$arr[1,11] = [double]"1"
$arr[1,12] = [double]"5.994110713235593e-15"
$arr[1,13] = [double]"6.660187513674559e-16"
$arr[1,14] = [double]"-4.043647085677147e-16"
$arr[1,15] = [double]"-9.656940896947432e-12"
$arr[2,0] = [double]"9.117503425189895e-18"
$arr[2,1] = [double]"-1.29335736330914e-23"
$arr[2,2] = [double]"-7.187995025634125e-12"
$arr[2,3] = [double]"7.188119076300633e-12"
$arr[2,4] = [double]"-1.11106675629875e-21"
$arr[2,5] = [double]"-8.729270386699821e-22"
$arr[2,6] = [double]"4.928775308882372e-22"
$arr[2,7] = [double]"2.443626581453735e-17"
$arr[2,8] = [double]"-1.50796173660397e-11"
$arr[2,9] = [double]"1.564309955763869e-11"
$arr[2,10] = [double]"2.414294852692855e-11"
$arr[2,11] = [double]"5.941540696923834e-17"
$arr[2,12] = [double]"6.106718605282241e-16"
$arr[2,13] = [double]"1.339266608542708e-15"
$arr[2,14] = [double]"1"
$arr[2,15] = [double]"-1.136173316801328e-11"
$arr[3,0] = [double]"-1.110568744945022e-16"
$arr[3,1] = [double]"2.931721386475353e-17"
$arr[3,2] = [double]"-1.517449074881796e-11"
$arr[3,3] = [double]"-1.51748937261747e-11"
$arr[3,4] = [double]"-2.414323856242113e-11"
$arr[3,5] = [double]"-1.211193598966174e-16"
$arr[3,6] = [double]"1.211050210022838e-16"
$arr[3,7] = [double]"2.195917947110037e-16"
$arr[3,8] = [double]"-0.7807480474653604"
$arr[3,9] = [double]"0.8099215765986721"
$arr[3,10] = [double]"6.122626597949276e-17"
$arr[3,11] = [double]"-2.759176053385392e-11"
$arr[3,12] = [double]"-1.574023087392656e-16"
$arr[3,13] = [double]"-8.995716286418389e-17"
$arr[3,14] = [double]"-1.931424581946772e-11"
$arr[3,15] = [double]"1.503558541149491e-22"
$arr[4,0] = [double]"-9.657188836838101e-12"
$arr[4,1] = [double]"-1.756118195063036e-21"
$arr[4,2] = [double]"1.319085049485686e-17"
$arr[4,3] = [double]"-1.319118218817301e-17"
$arr[4,4] = [double]"1.683647746509007e-21"
$arr[4,5] = [double]"7.799659062309739e-22"
$arr[4,6] = [double]"-3.572209909903872e-22"
$arr[4,7] = [double]"-1.121540905332171e-17"
$arr[4,8] = [double]"-1.34094410500141e-11"
$arr[4,9] = [double]"-1.258732465393586e-11"
$arr[4,10] = [double]"2.759189352235694e-11"
$arr[4,11] = [double]"-6.213372908313923e-16"
$arr[4,12] = [double]"-4.089820745680538e-15"
$arr[4,13] = [double]"1"
$arr[4,14] = [double]"-3.056582764667324e-15"
$arr[4,15] = [double]"-1.07286927637977e-11"
$arr[5,0] = [double]"-1"
$arr[5,1] = [double]"-1.931434256379727e-10"
$arr[5,2] = [double]"-2.79145909033353e-16"
$arr[5,3] = [double]"2.830006589978719e-16"
$arr[5,4] = [double]"6.43811430059132e-11"
$arr[5,5] = [double]"-1.51324093306471e-18"
$arr[5,6] = [double]"1.530787023740159e-18"
$arr[5,7] = [double]"1.401137829263262e-16"
$arr[5,8] = [double]"-2.554728460056816e-16"
$arr[5,9] = [double]"5.35963860758549e-17"
$arr[5,10] = [double]"7.722582432092188e-17"
$arr[5,11] = [double]"-1.073019029903051e-11"
$arr[5,12] = [double]"3.760756258047686e-23"
$arr[5,13] = [double]"-9.6571712691485e-12"
$arr[5,14] = [double]"1.022041877917098e-23"
$arr[5,15] = [double]"1.036063555944485e-22"
$arr[6,0] = [double]"1.532980422344469e-16"
$arr[6,1] = [double]"-7.593513408075343e-19"
$arr[6,2] = [double]"-9.151184712442485e-17"
$arr[6,3] = [double]"9.327291057562296e-17"
$arr[6,4] = [double]"-1.931440674250336e-11"
$arr[6,5] = [double]"1.517480723843754e-11"
$arr[6,6] = [double]"-1.517475196535439e-11"
$arr[6,7] = [double]"-3.951529720224298e-17"
$arr[6,8] = [double]"-1.085148659835906e-16"
$arr[6,9] = [double]"2.057559085583267e-16"
$arr[6,10] = [double]"1"
$arr[6,11] = [double]"-7.034652372932426e-19"
$arr[6,12] = [double]"4.285286258340713e-19"
$arr[6,13] = [double]"-2.759191801217852e-11"
$arr[6,14] = [double]"-2.414292816397255e-11"
$arr[6,15] = [double]"1.706409072799458e-22"
$arr[7,0] = [double]"-6.438079788469179e-11"
$arr[7,1] = [double]"-1.541994180325995e-15"
$arr[7,2] = [double]"-5.441272745586317e-16"
$arr[7,3] = [double]"1.864211858351288e-15"
$arr[7,4] = [double]"-1"
$arr[7,5] = [double]"-1.365703965208477e-10"
$arr[7,6] = [double]"-1.365756551291545e-10"
$arr[7,7] = [double]"6.51289742165428e-17"
$arr[7,8] = [double]"1.884941967829163e-11"
$arr[7,9] = [double]"-1.955375358015932e-11"
$arr[7,10] = [double]"-1.931427712973458e-11"
$arr[7,11] = [double]"4.933582038967252e-22"
$arr[7,12] = [double]"1.481658855806176e-25"
$arr[7,13] = [double]"4.231201798904439e-22"
$arr[7,14] = [double]"4.663175341856139e-22"
$arr[7,15] = [double]"-2.342812185182386e-33"
$arr[8,0] = [double]"7.782943328115727e-18"
$arr[8,1] = [double]"4.831206747893899e-22"
$arr[8,2] = [double]"-8.033655272116433e-12"
$arr[8,3] = [double]"8.033793917072827e-12"
$arr[8,4] = [double]"-4.089991652730208e-22"
$arr[8,5] = [double]"-5.420738390863429e-22"
$arr[8,6] = [double]"8.367266997640653e-22"
$arr[8,7] = [double]"1.931429166024468e-11"
$arr[8,8] = [double]"-1.508562958748424e-11"
$arr[8,9] = [double]"-1.416077417627182e-11"
$arr[8,10] = [double]"-1.228362757152175e-17"
$arr[8,11] = [double]"-5.513015296897995e-15"
$arr[8,12] = [double]"1"
$arr[8,13] = [double]"5.567700535660178e-15"
$arr[8,14] = [double]"-1.198007077412873e-15"
$arr[8,15] = [double]"-1.016482316582114e-11"
$arr[9,0] = [double]"7.422668723844378e-17"
$arr[9,1] = [double]"2.414292984068295e-11"
$arr[9,2] = [double]"-1.951051822093606e-11"
$arr[9,3] = [double]"-1.951033888895094e-11"
$arr[9,4] = [double]"7.132816766766335e-17"
$arr[9,5] = [double]"-1.274664482070851e-18"
$arr[9,6] = [double]"1.266471142440308e-18"
$arr[9,7] = [double]"1"
$arr[9,8] = [double]"8.586242091546045e-17"
$arr[9,9] = [double]"7.927753658747034e-17"
$arr[9,10] = [double]"-1.47451497343439e-17"
$arr[9,11] = [double]"-2.146024322006007e-11"
$arr[9,12] = [double]"-1.931443291938531e-11"
$arr[9,13] = [double]"2.24637446287797e-17"
$arr[9,14] = [double]"-2.065563270140633e-18"
$arr[9,15] = [double]"1.320040313855912e-22"
$arr[10,0] = [double]"1.600739220577689e-16"
$arr[10,1] = [double]"-2.522237112090684e-15"
$arr[10,2] = [double]"-0.7071006793750427"
$arr[10,3] = [double]"0.7071128825524579"
$arr[10,4] = [double]"7.044424980886948e-17"
$arr[10,5] = [double]"-6.828663451425932e-11"
$arr[10,6] = [double]"6.828638575721117e-11"
$arr[10,7] = [double]"-4.906190289695546e-16"
$arr[10,8] = [double]"-2.300957280430481e-16"
$arr[10,9] = [double]"3.488532565669355e-17"
$arr[10,10] = [double]"-1.937106113283201e-17"
$arr[10,11] = [double]"2.106310406990328e-24"
$arr[10,12] = [double]"-1.136137796367196e-11"
$arr[10,13] = [double]"8.442281831172406e-24"
$arr[10,14] = [double]"-1.016544344117642e-11"
$arr[10,15] = [double]"1.154900087443605e-22"
$arr[11,0] = [double]"-2.004483291426057e-25"
$arr[11,1] = [double]"4.734144518475036e-15"
$arr[11,2] = [double]"-0.7071128829453985"
$arr[11,3] = [double]"-0.7071006797679902"
$arr[11,4] = [double]"-1.923134077490348e-16"
$arr[11,5] = [double]"-6.828645794206624e-11"
$arr[11,6] = [double]"-6.828656788795495e-11"
$arr[11,7] = [double]"-2.759191791175646e-11"
$arr[11,8] = [double]"1.675503192581534e-11"
$arr[11,9] = [double]"-1.738122816398455e-11"
$arr[11,10] = [double]"-2.182431025882176e-25"
$arr[11,11] = [double]"5.921275742077835e-22"
$arr[11,12] = [double]"3.485224005666194e-22"
$arr[11,13] = [double]"4.35120607344181e-28"
$arr[11,14] = [double]"2.996742145423013e-22"
$arr[11,15] = [double]"-2.13436292205868e-33"
$arr[12,0] = [double]"6.37286102559884e-17"
$arr[12,1] = [double]"1.931421085411896e-11"
$arr[12,2] = [double]"3.562859036743195e-16"
$arr[12,3] = [double]"1.720734678578681e-17"
$arr[12,4] = [double]"-1.705925861010295e-16"
$arr[12,5] = [double]"1.951046815470161e-11"
$arr[12,6] = [double]"-1.951039708355994e-11"
$arr[12,7] = [double]"1.544619257860951e-16"
$arr[12,8] = [double]"-0.6248459701230595"
$arr[12,9] = [double]"-0.5865381826956547"
$arr[12,10] = [double]"1.672549891932644e-16"
$arr[12,11] = [double]"-8.869756629800233e-25"
$arr[12,12] = [double]"-2.414292817278483e-11"
$arr[12,13] = [double]"-2.14603805980339e-11"
$arr[12,14] = [double]"2.579890704434267e-25"
$arr[12,15] = [double]"1.49428971361494e-22"
$arr[13,0] = [double]"-1.931434253822939e-10"
$arr[13,1] = [double]"1"
$arr[13,2] = [double]"9.226093206815878e-16"
$arr[13,3] = [double]"4.078367769298119e-15"
$arr[13,4] = [double]"-1.12611528290278e-15"
$arr[13,5] = [double]"-4.552476976591992e-11"
$arr[13,6] = [double]"-4.552391412263053e-11"
$arr[13,7] = [double]"-2.414292817278628e-11"
$arr[13,8] = [double]"1.206842813504288e-11"
$arr[13,9] = [double]"1.132854215111745e-11"
$arr[13,10] = [double]"2.844064444401423e-25"
$arr[13,11] = [double]"3.657259261213452e-22"
$arr[13,12] = [double]"4.663069666697254e-22"
$arr[13,13] = [double]"3.163231711937532e-22"
$arr[13,14] = [double]"4.986840889354011e-29"
$arr[13,15] = [double]"-2.009515843114686e-33"
$arr[14,0] = [double]"-3.972885375306269e-26"
$arr[14,1] = [double]"5.834818475341819e-16"
$arr[14,2] = [double]"-6.828561017748119e-11"
$arr[14,3] = [double]"6.82874156034956e-11"
$arr[14,4] = [double]"2.591562785905889e-15"
$arr[14,5] = [double]"0.7071080691215335"
$arr[14,6] = [double]"-0.707105493245123"
$arr[14,7] = [double]"-6.285102656045612e-26"
$arr[14,8] = [double]"1.72406981169169e-11"
$arr[14,9] = [double]"1.618371951295341e-11"
$arr[14,10] = [double]"-2.146038059803278e-11"
$arr[14,11] = [double]"1.509657307716026e-29"
$arr[14,12] = [double]"4.571615443381234e-22"
$arr[14,13] = [double]"5.921330619709791e-22"
$arr[14,14] = [double]"3.897903770349174e-22"
$arr[14,15] = [double]"-2.473721641472796e-33"
$arr[15,0] = [double]"-1.243477728755538e-20"
$arr[15,1] = [double]"6.438114179409845e-11"
$arr[15,2] = [double]"-6.828710217201182e-11"
$arr[15,3] = [double]"-6.828592369026878e-11"
$arr[15,4] = [double]"-1.931434253822891e-10"
$arr[15,5] = [double]"0.7071054932492159"
$arr[15,6] = [double]"0.7071080691256262"
$arr[15,7] = [double]"1.998447998379254e-21"
$arr[15,8] = [double]"-1.800187511760553e-21"
$arr[15,9] = [double]"3.526542429072731e-22"
$arr[15,10] = [double]"8.750411324752787e-22"
$arr[15,11] = [double]"-2.002119097619507e-32"
$arr[15,12] = [double]"-1.637943704348285e-32"
$arr[15,13] = [double]"-1.607478975179168e-32"
$arr[15,14] = [double]"-1.313109907654873e-32"
$arr[15,15] = [double]"5.088256987469149e-44"
$ws.Range("A2:P17").Value = $arr
Write-Host "done"
